$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated IPC index values for B58:B97
$ws.Range("B58").Value = 99.24113860008084
$ws.Range("B59").Value = 100.08833502656637
$ws.Range("B60").Value = 100.9202492523565
$ws.Range("B61").Value = 101.3421504205054
$ws.Range("B62").Value = 102.1363794759627
$ws.Range("B63").Value = 103.3953126335217
$ws.Range("B64").Value = 104.06522469192201
$ws.Range("B65").Value = 105.40904228336889
$ws.Range("B66").Value = 106.37319359077487
$ws.Range("B67").Value = 106.59702069688022
$ws.Range("B68").Value = 106.81304646207194
$ws.Range("B69").Value = 107.70379628281218
$ws.Range("B70").Value = 107.93621988160507
$ws.Range("B71").Value = 107.95448866196907
$ws.Range("B72").Value = 109.06297792005627
$ws.Range("B73").Value = 110.18847889914778
$ws.Range("B74").Value = 111.28573001376111
$ws.Range("B75").Value = 113.05489218575623
$ws.Range("B76").Value = 114.04726496638703
$ws.Range("B77").Value = 115.8595659225107
$ws.Range("B78").Value = 116.15342921506523
$ws.Range("B79").Value = 117.77539459864363
$ws.Range("B80").Value = 119.46307734324604
$ws.Range("B81").Value = 120.61112184976389
$ws.Range("B82").Value = 121.55465333769331
$ws.Range("B83").Value = 122.73896399367251
$ws.Range("B84").Value = 123.63667132280123
$ws.Range("B85").Value = 124.02259631937034
$ws.Range("B86").Value = 124.90266543904843
$ws.Range("B87").Value = 125.56318686353885
$ws.Range("B88").Value = 125.69747726963409
$ws.Range("B89").Value = 126.52672985443553
$ws.Range("B90").Value = 127.39592550671495
$ws.Range("B91").Value = 128.27415647484284
$ws.Range("B92").Value = 129.23121015536356
$ws.Range("B93").Value = 130.02431410600929
$ws.Range("B94").Value = 129.94174979755198
$ws.Range("B95").Value = 131.25481285276797
$ws.Range("B96").Value = 132.027202885248
$ws.Range("B97").Value = 133.31451960604798

# Remove the right+bottom border previously applied to the data column (B2:B97)
$ws.Range("B2:B97").Borders.LineStyle = -4142

